$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.56"
$ws.Range("E2").Value = "'6.05%"

$ws.Range("D3").Value = "'32.21"
$ws.Range("E3").Value = "'9.81%"

$ws.Range("D4").Value = "'5.318"
$ws.Range("E4").Value = "'4.25%"

$ws.Range("D5").Value = "'0.07415"
$ws.Range("E5").Value = "'11.03%"

$ws.Range("D6").Value = "'7.756"
$ws.Range("E6").Value = "'5.33%"

$ws.Range("D7").Value = "'3.687"
$ws.Range("E7").Value = "'8.09%"

$ws.Range("D8").Value = "'1.529"
$ws.Range("E8").Value = "'13.18%"

$ws.Range("D9").Value = "'0.9079"
$ws.Range("E9").Value = "'-0.97%"

$ws.Range("D10").Value = "'0.01656"
$ws.Range("E10").Value = "'2,471.68%"

$ws.Range("D11").Value = "'0.1669"
$ws.Range("E11").Value = "'5.10%"

$ws.Range("D12").Value = "'0.07610"
$ws.Range("E12").Value = "'13.34%"

$ws.Range("D13").Value = "'0.08005"
$ws.Range("E13").Value = "'4.20%"

$ws.Range("D14").Value = "'0.03019"
$ws.Range("E14").Value = "'2.50%"

$ws.Range("D15").Value = "'0.09842"
$ws.Range("E15").Value = "'9.48%"

$ws.Range("D16").Value = "'0.001521"
$ws.Range("E16").Value = "'-3.16%"

$ws.Range("D17").Value = "'0.04561"
$ws.Range("E17").Value = "'0.97%"

$ws.Range("D18").Value = "'0.006322"
$ws.Range("E18").Value = "'1.07%"

$ws.Range("D19").Value = "'3.495"
$ws.Range("E19").Value = "'1.39%"

$ws.Range("D20").Value = "'2.240"

$ws.Range("D21").Value = "'0.3264"
$ws.Range("E21").Value = "'1.60%"

$ws.Range("D22").Value = "'0.1332"
$ws.Range("E22").Value = "'1.73%"

$ws.Range("D23").Value = "'4.215"
$ws.Range("E23").Value = "'3.47%"

$ws.Range("D25").Value = "'0.001213"
$ws.Range("E25").Value = "'1.88%"

$ws.Range("D26").Value = "'0.004507"
$ws.Range("E26").Value = "'9.23%"

$ws.Range("D27").Value = "'0.0001169"
$ws.Range("E27").Value = "'-6.48%"

$ws.Range("D28").Value = "'0.0001736"
$ws.Range("E28").Value = "'7.27%"

$ws.Range("D40").Value = "'0.04499"
$ws.Range("E40").Value = "'6.86%"

$ws.Range("D41").Value = "'0.007286"
$ws.Range("E41").Value = "'8.23%"

$ws.Range("E42").Value = "'9.62%"

$ws.Range("D43").Value = "'0.002258"
$ws.Range("E43").Value = "'14.08%"

$ws.Range("E44").Value = "'1.37%"

$ws.Range("D45").Value = "'0.00006124"
$ws.Range("E45").Value = "'7.50%"

$ws.Range("D46").Value = "'1.893"
$ws.Range("E46").Value = "'-3.94%"

$ws.Range("D47").Value = "'0.01296"
$ws.Range("E47").Value = "'-0.84%"
